$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.99 -> 0M
$t.Cell(1, 1).Range.Text = "0M"

# Row 2: 0 -> 0M
$t.Cell(2, 1).Range.Text = "0M"

# Row 3: 75 -> 0M
$t.Cell(3, 1).Range.Text = "0M"

# Insert 10 new rows right after row 3 (before the old row 4), each with a
# single value, in order. Each Rows.Add(refRow) call inserts immediately
# before refRow, so walking the desired values backwards (and always
# inserting before the same fixed row) reproduces the original order.
$newValues = @("105", "0.00002", "0.00006", "0.00004", "0.00001", "0.00004", "0.00004", "0.00006", "0.00442", "100.0")
$refRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
}

# The three rows that used to hold 9 tab-separated numbers collapse down to a
# single short value each. After the 10-row insertion above, these are rows
# 44, 45, 46.
$t.Cell(44, 1).Range.Text = "99.99"
$t.Cell(45, 1).Range.Text = "0"
$t.Cell(46, 1).Range.Text = "75"
